# Mark all "ORD" (on-order) status entries in the BOM as "INV" (in inventory),
# matching the highlighted "INV" formatting already used elsewhere in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows in the Status column (B) that currently read "ORD" and must become "INV",
# with the same yellow highlight / center alignment already used by the other
# "INV" rows (e.g. row 16).
$rows = @(8,10,11,12,13,14,15,17,18,19,20,21,22,23,24,25,26,27,28,30,31)

foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r, 2)
    $cell.Value = "INV"
    $cell.Interior.ColorIndex = 6
    $cell.HorizontalAlignment = -4108
}

# Restore the cursor/selection to the cell left selected by the author.
$ws.Range("C37").Select()
